$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "DSD: HTS_TST (Facility)"
$ws.Range("E2").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A5 Devia ter o valor: 'DSD: TX_NEW'"

# Row 3
$ws.Range("D3").Value = "Subtotal"
$ws.Range("E3").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A21 Devia ter o valor: 'DSD: TX_CURR'"

# Row 4
$ws.Range("D4").Value = "Positive"
$ws.Range("E4").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A44 Devia ter o valor: 'DSD: TX_RTT'"

# Row 5
$ws.Range("D5").Value = "Subtotal"
$ws.Range("E5").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A62 Devia ter o valor: 'DSD: TX_ML'"

# Row 6
$ws.Range("D6").Value = "Positive"
$ws.Range("E6").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A99 Devia ter o valor: 'DSD: PMTCT_ART'"

# Row 7 - value_on_template cell is removed, only error_message updated
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A114 Devia ter o valor: 'DSD: TX_PVLS (Numerator)'"

# Row 8
$ws.Range("D8").Value = "Female"
$ws.Range("E8").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A134 Devia ter o valor: 'DSD:TX_PVLS (Denominator)'"

# Row 9
$ws.Range("D9").Value = "Auto-Calculate"
$ws.Range("E9").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A164 Devia ter o valor: 'DSD: TX_TB (Denominator)'"

# Row 10 - value_on_template cell is removed, only error_message updated
$ws.Range("D10").ClearContents()
$ws.Range("E10").Value = "2022-08-31 TEMPLATE ERROR ('MER C&T|MER_ATS_Xipamanine_12'): O ficheiro de importacao nao esta consistente, a cellula: A182 Devia ter o valor: 'DSD: TB_ART (Numerator)'"
